$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.212.46"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "1.896.18"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  +0.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5199"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07284"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("E11").Value = "  +1.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08185"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.67%  "

$ws.Range("D14").Value = "1.896.64"
$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.281"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008615"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "

$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").Value = "27.243.19"
$ws.Range("E20").Value = "  +0.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.086"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.403"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.295"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.746"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.72%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.952"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("E30").Value = "  +0.83%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09227"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.41%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05032"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7949"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.25%  "

$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.445"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.946"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.595"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5664"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01986"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.953"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.558"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1516"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4894"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.624"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.69%  "

$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05943"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
